# Auto-generated: applies updated TPM-derived values to F11r-Jam2 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 27.636609
$ws.Range("H2").Value = 82.90982700000001
$ws.Range("I2").Value = 0.5611737787305786
$ws.Range("J2").Value = 0.5623609389937425
$ws.Range("M2").Value = 54.23134333333334
$ws.Range("N2").Value = 162.69403
$ws.Range("O2").Value = 0.9097185042023884
$ws.Range("P2").Value = 0.9200351849746305
$ws.Range("Q2").Value = 1498.77043124809
$ws.Range("R2").Value = 13488.93388123281
$ws.Range("S2").Value = 0.5105101705843841
$ws.Range("T2").Value = 0.5173918505296148

# Row 3
$ws.Range("G3").Value = 27.636609
$ws.Range("H3").Value = 82.90982700000001
$ws.Range("I3").Value = 0.5611737787305786
$ws.Range("J3").Value = 0.5623609389937425
$ws.Range("O3").Value = 0.04961247197704183
$ws.Range("P3").Value = 0.0501751032012552
$ws.Range("Q3").Value = 81.73704907267901
$ws.Range("R3").Value = 735.633441654111
$ws.Range("S3").Value = 0.0278412183715215
$ws.Range("T3").Value = 0.02821651815036581

# Row 4
$ws.Range("G4").Value = 27.636609
$ws.Range("H4").Value = 82.90982700000001
$ws.Range("I4").Value = 0.5611737787305786
$ws.Range("J4").Value = 0.5623609389937425
$ws.Range("M4").Value = 0.2447093333333333
$ws.Range("N4").Value = 0.734128
$ws.Range("O4").Value = 0.004104943654374356
$ws.Range("P4").Value = 0.004151495849448536
$ws.Range("Q4").Value = 6.762936163984001
$ws.Range("R4").Value = 60.866425475856
$ws.Range("S4").Value = 0.002303586742001367
$ws.Range("T4").Value = 0.002334639104124503

# Row 5
$ws.Range("G5").Value = 27.636609
$ws.Range("H5").Value = 82.90982700000001
$ws.Range("I5").Value = 0.5611737787305786
$ws.Range("J5").Value = 0.5623609389937425
$ws.Range("M5").Value = 2.005396
$ws.Range("N5").Value = 4.010792
$ws.Range("O5").Value = 0.03364006379558217
$ws.Range("P5").Value = 0.02268103973830366
$ws.Range("Q5").Value = 55.42234514216401
$ws.Range("R5").Value = 332.534070852984
$ws.Range("S5").Value = 0.01887792171690458
$ws.Range("T5").Value = 0.01275493080458684

# Row 6
$ws.Range("G6").Value = 27.636609
$ws.Range("H6").Value = 82.90982700000001
$ws.Range("I6").Value = 0.5611737787305786
$ws.Range("J6").Value = 0.5623609389937425
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1743103333333333
$ws.Range("N6").Value = 0.522931
$ws.Range("O6").Value = 0.002924016370613348
$ws.Range("P6").Value = 0.002957176236362014
$ws.Range("Q6").Value = 4.817346526993001
$ws.Range("R6").Value = 43.35611874293701
$ws.Range("S6").Value = 0.001640881315767165
$ws.Range("T6").Value = 0.001663000405050523

# Row 7
$ws.Range("G7").Value = 11.07680033333333
$ws.Range("H7").Value = 33.230401
$ws.Range("I7").Value = 0.224919413930298
$ws.Range("J7").Value = 0.2253952298018738
$ws.Range("M7").Value = 54.23134333333334
$ws.Range("N7").Value = 162.69403
$ws.Range("O7").Value = 0.9097185042023884
$ws.Range("P7").Value = 0.9200351849746305
$ws.Range("Q7").Value = 600.7097619117811
$ws.Range("R7").Value = 5406.38785720603
$ws.Range("S7").Value = 0.2046133528067485
$ws.Range("T7").Value = 0.2073715419431663

# Row 8
$ws.Range("G8").Value = 11.07680033333333
$ws.Range("H8").Value = 33.230401
$ws.Range("I8").Value = 0.224919413930298
$ws.Range("J8").Value = 0.2253952298018738
$ws.Range("O8").Value = 0.04961247197704183
$ws.Range("P8").Value = 0.0501751032012552
$ws.Range("Q8").Value = 32.76034959332144
$ws.Range("R8").Value = 294.843146339893
$ws.Range("S8").Value = 0.01115880812070958
$ws.Range("T8").Value = 0.01130922891637965

# Row 9
$ws.Range("G9").Value = 11.07680033333333
$ws.Range("H9").Value = 33.230401
$ws.Range("I9").Value = 0.224919413930298
$ws.Range("J9").Value = 0.2253952298018738
$ws.Range("M9").Value = 0.2447093333333333
$ws.Range("N9").Value = 0.734128
$ws.Range("O9").Value = 0.004104943654374356
$ws.Range("P9").Value = 0.004151495849448536
$ws.Range("Q9").Value = 2.710596425036444
$ws.Range("R9").Value = 24.395367825328
$ws.Range("S9").Value = 0.0009232815209587758
$ws.Range("T9").Value = 0.000935727361007978

# Row 10
$ws.Range("G10").Value = 11.07680033333333
$ws.Range("H10").Value = 33.230401
$ws.Range("I10").Value = 0.224919413930298
$ws.Range("J10").Value = 0.2253952298018738
$ws.Range("M10").Value = 2.005396
$ws.Range("N10").Value = 4.010792
$ws.Range("O10").Value = 0.03364006379558217
$ws.Range("P10").Value = 0.02268103973830366
$ws.Range("Q10").Value = 22.21337108126534
$ws.Range("R10").Value = 133.280226487592
$ws.Range("S10").Value = 0.007566303433480177
$ws.Range("T10").Value = 0.005112198163960386

# Row 11
$ws.Range("G11").Value = 11.07680033333333
$ws.Range("H11").Value = 33.230401
$ws.Range("I11").Value = 0.224919413930298
$ws.Range("J11").Value = 0.2253952298018738
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1743103333333333
$ws.Range("N11").Value = 0.522931
$ws.Range("O11").Value = 0.002924016370613348
$ws.Range("P11").Value = 0.002957176236362014
$ws.Range("Q11").Value = 1.930800758370111
$ws.Range("R11").Value = 17.377206825331
$ws.Range("S11").Value = 0.0006576680484009514
$ws.Range("T11").Value = 0.0006665334173594563

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.311891
$ws.Range("H12").Value = 0.6237819999999999
$ws.Range("I12").Value = 0.006333087066581101
$ws.Range("J12").Value = 0.00423098978661956
$ws.Range("M12").Value = 54.23134333333334
$ws.Range("N12").Value = 162.69403
$ws.Range("O12").Value = 0.9097185042023884
$ws.Range("P12").Value = 0.9200351849746305
$ws.Range("Q12").Value = 16.91426790357666
$ws.Range("R12").Value = 101.48560742146
$ws.Range("S12").Value = 0.005761326493193651
$ws.Range("T12").Value = 0.0038926594709583

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.311891
$ws.Range("H13").Value = 0.6237819999999999
$ws.Range("I13").Value = 0.006333087066581101
$ws.Range("J13").Value = 0.00423098978661956
$ws.Range("O13").Value = 0.04961247197704183
$ws.Range("P13").Value = 0.0501751032012552
$ws.Range("Q13").Value = 0.9224376974876666
$ws.Range("R13").Value = 5.534626184925999
$ws.Range("S13").Value = 0.0003142001046189209
$ws.Range("T13").Value = 0.0002122903491870932

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.311891
$ws.Range("H14").Value = 0.6237819999999999
$ws.Range("I14").Value = 0.006333087066581101
$ws.Range("J14").Value = 0.00423098978661956
$ws.Range("M14").Value = 0.2447093333333333
$ws.Range("N14").Value = 0.734128
$ws.Range("O14").Value = 0.004104943654374356
$ws.Range("P14").Value = 0.004151495849448536
$ws.Range("Q14").Value = 0.07632263868266666
$ws.Range("R14").Value = 0.457935832096
$ws.Range("S14").Value = 0.00002599696556656239
$ws.Range("T14").Value = 0.00001756493653821025

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.311891
$ws.Range("H15").Value = 0.6237819999999999
$ws.Range("I15").Value = 0.006333087066581101
$ws.Range("J15").Value = 0.00423098978661956
$ws.Range("M15").Value = 2.005396
$ws.Range("N15").Value = 4.010792
$ws.Range("O15").Value = 0.03364006379558217
$ws.Range("P15").Value = 0.02268103973830366
$ws.Range("Q15").Value = 0.625464963836
$ws.Range("R15").Value = 2.501859855344
$ws.Range("S15").Value = 0.0002130454529427646
$ws.Range("T15").Value = 0.0000959632474826752

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.311891
$ws.Range("H16").Value = 0.6237819999999999
$ws.Range("I16").Value = 0.006333087066581101
$ws.Range("J16").Value = 0.00423098978661956
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1743103333333333
$ws.Range("N16").Value = 0.522931
$ws.Range("O16").Value = 0.002924016370613348
$ws.Range("P16").Value = 0.002957176236362014
$ws.Range("Q16").Value = 0.05436582417366666
$ws.Range("R16").Value = 0.326194945042
$ws.Range("S16").Value = 0.00001851805025920281
$ws.Range("T16").Value = 0.00001251178245328175

# Row 17
$ws.Range("G17").Value = 10.222562
$ws.Range("H17").Value = 30.667686
$ws.Range("I17").Value = 0.2075737202725422
$ws.Range("J17").Value = 0.208012841417764
$ws.Range("M17").Value = 54.23134333333334
$ws.Range("N17").Value = 162.69403
$ws.Range("O17").Value = 0.9097185042023884
$ws.Range("P17").Value = 0.9200351849746305
$ws.Range("Q17").Value = 554.3832695682868
$ws.Range("R17").Value = 4989.44942611458
$ws.Range("S17").Value = 0.1888336543180621
$ws.Range("T17").Value = 0.191379133030891

# Row 18
$ws.Range("G18").Value = 10.222562
$ws.Range("H18").Value = 30.667686
$ws.Range("I18").Value = 0.2075737202725422
$ws.Range("J18").Value = 0.208012841417764
$ws.Range("O18").Value = 0.04961247197704183
$ws.Range("P18").Value = 0.0501751032012552
$ws.Range("Q18").Value = 30.23388476648867
$ws.Range("R18").Value = 272.104962898398
$ws.Range("S18").Value = 0.01029824538019182
$ws.Range("T18").Value = 0.01043706578532264

# Row 19
$ws.Range("G19").Value = 10.222562
$ws.Range("H19").Value = 30.667686
$ws.Range("I19").Value = 0.2075737202725422
$ws.Range("J19").Value = 0.208012841417764
$ws.Range("M19").Value = 0.2447093333333333
$ws.Range("N19").Value = 0.734128
$ws.Range("O19").Value = 0.004104943654374356
$ws.Range("P19").Value = 0.004151495849448536
$ws.Range("Q19").Value = 2.501556331978667
$ws.Range("R19").Value = 22.514006987808
$ws.Range("S19").Value = 0.0008520784258476495
$ws.Range("T19").Value = 0.0008635644477778439

# Row 20
$ws.Range("G20").Value = 10.222562
$ws.Range("H20").Value = 30.667686
$ws.Range("I20").Value = 0.2075737202725422
$ws.Range("J20").Value = 0.208012841417764
$ws.Range("M20").Value = 2.005396
$ws.Range("N20").Value = 4.010792
$ws.Range("O20").Value = 0.03364006379558217
$ws.Range("P20").Value = 0.02268103973830366
$ws.Range("Q20").Value = 20.50028494455201
$ws.Range("R20").Value = 123.001709667312
$ws.Range("S20").Value = 0.006982793192254646
$ws.Range("T20").Value = 0.004717947522273765

# Row 21
$ws.Range("G21").Value = 10.222562
$ws.Range("H21").Value = 30.667686
$ws.Range("I21").Value = 0.2075737202725422
$ws.Range("J21").Value = 0.2253952298018738
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1743103333333333
$ws.Range("N21").Value = 0.522931
$ws.Range("O21").Value = 0.002924016370613348
$ws.Range("P21").Value = 0.002957176236362014
$ws.Range("Q21").Value = 1.781898189740667
$ws.Range("R21").Value = 16.037083707666
$ws.Range("S21").Value = 0.0006069489561860292
$ws.Range("T21").Value = 0.0006151306314987519

Write-Host "Applied 250 cell updates"
